$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row 2 (pushes the existing row 2 "Contact" and row 3 "About Me" down)
$ws.Rows.Item(2).Insert()

# Populate the newly inserted row with the SublimeText Package article entry
$ws.Range("A2").Value = "tag-web"
$ws.Range("B2").Value = "SublimeText Package"
$ws.Range("C2").Value = "/articles/201405/sublimetext-package.html"
$ws.Range("D2").Value = "/img/articles/201405/20140520_1_01.jpg"
$ws.Range("E2").Value = "MAY 20TH, 2014"

# Copy the "prettify" style (border + text format) from the header/body rows onto the new row
$ws.Range("A1:E1").Copy()
$ws.Range("A2:E2").PasteSpecial(-4122)

# Match the recorded selection from the edit
$ws.Range("E2").Select()
